$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I4").Value = 5
$ws.Range("Q4").Value = 2.08
$ws.Range("F6").Value = 2.7
$ws.Range("G6").Value = 4.1
$ws.Range("H6").Value = 2.24
$ws.Range("J6").Value = 2.44
$ws.Range("K6").Value = 5.6
$ws.Range("F7").Value = 2.12
$ws.Range("G7").Value = 2.4
$ws.Range("H7").Value = 3.4
$ws.Range("I7").Value = 4.1
$ws.Range("J7").Value = 3
$ws.Range("K7").Value = 4.1
$ws.Range("P7").Value = 1.81
$ws.Range("Q7").Value = 1.96
$ws.Range("G8").Value = 7.6
$ws.Range("P8").Value = 2.1
$ws.Range("F9").Value = 1.92
$ws.Range("H9").Value = 3.35
$ws.Range("P9").Value = 2.22
$ws.Range("Q9").Value = 1.5
$ws.Range("F10").Value = 5.3
$ws.Range("H10").Value = 1.69
$ws.Range("I10").Value = 1.7
$ws.Range("K10").Value = 4.5
$ws.Range("O10").Value = 1.24
$ws.Range("P10").Value = 2.34
$ws.Range("V10").Value = 2.42
$ws.Range("U11").Value = 2.12
$ws.Range("Y11").Value = 14
$ws.Range("AD11").Value = 16
$ws.Range("AE11").Value = 55
$ws.Range("AF11").Value = 14
$ws.Range("AJ11").Value = 36
$ws.Range("AK11").Value = 24
$ws.Range("G12").Value = 1.87
$ws.Range("P12").Value = 2.96
$ws.Range("Q12").Value = 1.45
$ws.Range("X12").Value = 36
$ws.Range("Y12").Value = 42
$ws.Range("AA12").Value = 80
$ws.Range("AB12").Value = 17
$ws.Range("AC12").Value = 11.5
$ws.Range("AD12").Value = 22
$ws.Range("AE12").Value = 980
$ws.Range("AF12").Value = 16.5
$ws.Range("AK12").Value = 17
$ws.Range("AN12").Value = 6.8
$ws.Range("AO12").Value = 42
$ws.Range("F13").Value = 1.76
$ws.Range("G13").Value = 1.78
$ws.Range("I13").Value = 5
$ws.Range("K13").Value = 4.4
$ws.Range("L13").Value = 1.33
$ws.Range("P13").Value = 2.22
$ws.Range("Q13").Value = 1.79
$ws.Range("S13").Value = 2.88
$ws.Range("AO13").Value = 1000
$ws.Range("F14").Value = 1.91
$ws.Range("G14").Value = 2.16
$ws.Range("I14").Value = 5.2
$ws.Range("M14").Value = 1.07
$ws.Range("N14").Value = 3.2
$ws.Range("Q14").Value = 2.06
$ws.Range("R14").Value = 1.29
$ws.Range("S14").Value = 3.4
$ws.Range("T14").Value = 1.86
$ws.Range("U14").Value = 1.92
$ws.Range("V14").Value = 1.25
$ws.Range("W14").Value = 1.87
$ws.Range("X14").Value = 980
$ws.Range("Y14").Value = 980
$ws.Range("Z14").Value = 980
$ws.Range("AA14").Value = 130
$ws.Range("AB14").Value = 10
$ws.Range("AC14").Value = 9.6
$ws.Range("AD14").Value = 980
$ws.Range("AE14").Value = 80
$ws.Range("AF14").Value = 980
$ws.Range("AG14").Value = 980
$ws.Range("AH14").Value = 980
$ws.Range("AI14").Value = 85
$ws.Range("AJ14").Value = 980
$ws.Range("AK14").Value = 980
$ws.Range("AL14").Value = 980
$ws.Range("AM14").Value = 150
$ws.Range("AN14").Value = 980
$ws.Range("AO14").Value = 95
$ws.Range("F16").Value = 1.25
$ws.Range("G16").Value = 980
$ws.Range("H16").Value = 1.25
$ws.Range("J16").Value = 1.09
$ws.Range("K16").Value = 4.9

